# Change "hastigt skiftende omstændigheder" -> "hurtigt skiftende omstændigheder"
# in the CV bullet about the trip to Southeast Asia.
#
# The run that holds "hastigt" also holds a lot of surrounding text
# (".. tilpasse mig hastigt skiftende ... kulturer. "). The target edit
# splits that single run into three runs: the text before "hastigt", the
# replacement word itself ("hurtigt"), and the text after it - while the
# earlier, unrelated run ("Rejsen udviklede min problemløsning") must stay
# exactly as it was (a separate run).

$d = $word.ActiveDocument

# Locate the word to replace.
$searchRange = $d.Content.Duplicate
$find = $searchRange.Find
$found = $find.Execute("hastigt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $wordStart = $searchRange.Start
    $wordEnd = $searchRange.End

    # Remember the extent of the preceding run ("Rejsen udviklede min
    # problemløsning") so it can be re-isolated after the edit below -
    # replacing text in the paragraph causes the engine to coalesce
    # adjacent same-formatted runs, which would otherwise merge this
    # untouched text into the edited run.
    $prefixSearchRange = $d.Content.Duplicate
    $prefixSearchRange.End = $wordStart
    $prefixFind = $prefixSearchRange.Find
    $prefixFound = $prefixFind.Execute("Rejsen udviklede min problemløsning", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # 1. Replace "hastigt" with "hurtigt".
    $d.Range($wordStart, $wordEnd).Text = "hurtigt"
    $newWordEnd = $wordStart + 7

    # 2. Re-isolate "hurtigt" into its own run, separate from the text
    #    immediately before and after it, by toggling (and restoring) a
    #    character property over just that word.
    $midRange = $d.Range($wordStart, $newWordEnd)
    $midRange.Bold = 1
    $midRange.Bold = 0

    # 3. Likewise re-isolate the earlier, untouched run so it does not stay
    #    merged with the text that now precedes "hurtigt".
    if ($prefixFound) {
        $prefixRange = $d.Range($prefixSearchRange.Start, $prefixSearchRange.End)
        $prefixRange.Bold = 1
        $prefixRange.Bold = 0
    }

    Write-Host "Replaced 'hastigt' with 'hurtigt' at $wordStart-$newWordEnd"
} else {
    Write-Host "WARNING: 'hastigt' not found; no changes made"
}
